$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 170, shifting existing rows 170:292 down to 171:293
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with the new record's data
$ws.Range("A170").Value = 10
$ws.Range("B170").Value = "Vega Modelo de Temuco"
$ws.Range("C170").Value = "La Araucanía"
$ws.Range("D170").Value = 45068
$ws.Range("E170").Value = 9
$ws.Range("F170").Value = 100112005
$ws.Range("G170").Value = "Puerro"
$ws.Range("H170").Value = "Azul de Maquehue"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 30
$ws.Range("K170").Value = 11000
$ws.Range("L170").Value = 11000
$ws.Range("M170").Value = 11000
$ws.Range("N170").Value = "$/docena de paquetes"
$ws.Range("O170").Value = "Provincia de Cautín"
$ws.Range("P170").Value = 917
$ws.Range("Q170").Value = 12
$ws.Range("R170").Value = "Hortaliza"
